$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new BOM rows. Done top-to-bottom using the row numbers of the
# *evolving* sheet (each subsequent Insert() targets positions that already
# account for the earlier inserts), so the final layout lines up with the
# target:
#   - rows 5:6  -> new CON1 entry + its trailing blank separator
#   - row 11    -> new C18 entry (lands at row 10 after the first insert
#                  shifts everything below row 4 down by two)
#   - rows 14:15-> new D3 / D4 entries (inserted right after D2)
$ws.Rows("5:6").Insert()
$ws.Rows("11:11").Insert()
$ws.Rows("14:15").Insert()

# --- Row 5: CON1 / DC power jack ---
$ws.Range("B5").Value = "CON1"
$ws.Range("C5").Value = "CONN PWR JACK 2.0X5.5MM SMD"
$ws.Range("E5").Value = "CP-070AHPJCT-ND"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "footprint:PJ-070AH-SMT"
$ws.Range("C5:G5").Style = "Normal"

# --- Row 10: C18 / tantalum cap ---
$ws.Range("B10").Value = "C18"
$ws.Range("C10").Value = "CAP TANT 100UF 10V 20% 2917"
$ws.Range("D10").Value = "DigiKey"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "399-3772-1-ND"
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = "Capacitors_Tantalum_SMD:TantalC_SizeD_EIA-7343_Reflow"

# --- Row 14: D3 / zener diode ---
$ws.Range("B14").Value = "D3"
$ws.Range("C14").Value = "DIODE ZENER 6.2V 500MW SOD123"
$ws.Range("D14").Value = "DigiKey"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "DDZ6V2BDICT-ND"
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = "Diodes_SMD:SOD-123"

# --- Row 15: D4 / general purpose diode ---
$ws.Range("B15").Value = "D4"
$ws.Range("C15").Value = "DIODE GEN PURP 100V 1A SMA"
$ws.Range("D15").Value = "DigiKey"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "S1B-FDICT-ND"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = "Diodes_SMD:SMA_Standard"

# Match the author's final selection/cursor position.
$ws.Range("B26").Select()
